$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "286.58"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.60%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "15"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29.37"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.84%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "15"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.075"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.51%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "15"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06770"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.33%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "15"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.367"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.94%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "15"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.375"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.71%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "15"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9000"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.89%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "15"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1589"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.97%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "15"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06896"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "11.83%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "15"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07648"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.69%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "15"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.17%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "15"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08994"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.25%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "15"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001572"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.44%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "15"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006438"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.72%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "15"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006254"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.49%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "15"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.458"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.47%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "15"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.440"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.82%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "15"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.229"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.03%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "15"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3200"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.32%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "15"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1322"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3.24%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "15"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.999"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.79%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "15"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1557"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.81%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "15"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04480"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "15"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001203"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.49%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "15"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004365"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.72%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "15"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001165"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-6.74%"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "15"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001768"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "9.32%"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "15"

$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "15"

$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "15"

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "15"

$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "15"

$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "15"

$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "15"

$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "15"

$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "15"

$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "15"

$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "15"

$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "15"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04273"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.77%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "15"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006779"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.89%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "15"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.77%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "15"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002191"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "7.45%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "15"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01150"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.73%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "15"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005716"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.07%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "15"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.87%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "15"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01302"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.15%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "15"

$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "15"

$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "15"

$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "15"

$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "15"
